$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell D1 ("function"), styled like B1/C1 (bold + border) ---
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D1").Value = "function"

# --- Data rows: centers/sigmas become text (shared strings) and a new
#     "function" column (D) with literal value "lorentzian" is added.
#     Values are entered as a formula producing a text literal, then
#     frozen to a static value via copy/paste-values, which keeps the
#     cell typed as text without altering its number format / style. ---
$centers = @(
    "589.0804626270364",
    "590.3215411282541",
    "587.0794956841607",
    "587.6397721021245",
    "584.6736011532843",
    "582.8119655993378"
)
$sigmas = @(
    "16.59078517060218",
    "22.04750045153666",
    "20.432815968077215",
    "19.522961422094134",
    "24.862970093436694",
    "25.064423076262877"
)

for ($i = 0; $i -lt 6; $i++) {
    $row = $i + 2

    $bCell = $ws.Cells.Item($row, 2)
    $bCell.Formula = "=""" + $centers[$i] + """"
    $bCell.Copy()
    $bCell.PasteSpecial(-4163)   # xlPasteValues

    $cCell = $ws.Cells.Item($row, 3)
    $cCell.Formula = "=""" + $sigmas[$i] + """"
    $cCell.Copy()
    $cCell.PasteSpecial(-4163)   # xlPasteValues

    $ws.Cells.Item($row, 4).Value = "lorentzian"
}
